# Apply crypto price/volume updates as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Several updated "Price" values are plain decimal-looking strings (e.g. "38.30",
# "591.59"). If such strings are assigned directly, Excel auto-converts them to
# floating point numbers, which can both drop a significant trailing zero and
# introduce binary floating point rounding noise. Force a text number format on
# those specific Price cells first so the exact original text is preserved.
$textCells = @("D4", "D5", "D6", "D11", "D12", "D13", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D29", "D33", "D34", "D35", "D37", "D38", "D39", "D40", "D41", "D42", "D45", "D46", "D47", "D49", "D50")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '67.353.56'
$ws.Range("E2").Value = '  -4.61%  '
$ws.Range("D3").Value = '3.265.38'
$ws.Range("E3").Value = '  -6.97%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '591.59'
$ws.Range("E5").Value = '  -5.12%  '
$ws.Range("D6").Value = '151.31'
$ws.Range("E6").Value = '  -11.79%  '
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").Value = '3.256.81'
$ws.Range("E8").Value = '  -7.06%  '
$ws.Range("E9").Value = '  -10.71%  '
$ws.Range("E10").Value = '  -14.16%  '
$ws.Range("D11").Value = '6.63'
$ws.Range("E11").Value = '  -7.71%  '
$ws.Range("D12").Value = '0.513'
$ws.Range("E12").Value = '  -12.16%  '
$ws.Range("D13").Value = '38.30'
$ws.Range("E13").Value = '  -16.92%  '
$ws.Range("E14").Value = '  -11.50%  '
$ws.Range("D15").Value = '3.784.50'
$ws.Range("E15").Value = '  -7.21%  '
$ws.Range("D16").Value = '67.305.90'
$ws.Range("D17").Value = '3.267.00'
$ws.Range("E17").Value = '  -7.11%  '
$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").Value = '7.24'
$ws.Range("E18").Value = '  -13.91%  '
$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").Value = '0.114'
$ws.Range("E19").Value = '  -6.18%  '
$ws.Range("D20").Value = '533.30'
$ws.Range("E20").Value = '  -11.86%  '
$ws.Range("D21").Value = '15.10'
$ws.Range("E21").Value = '  -14.48%  '
$ws.Range("D22").Value = '0.762'
$ws.Range("E22").Value = '  -13.17%  '
$ws.Range("D23").Value = '7.92'
$ws.Range("E23").Value = '  -13.04%  '
$ws.Range("D24").Value = '85.69'
$ws.Range("E24").Value = '  -11.69%  '
$ws.Range("D25").Value = '13.61'
$ws.Range("E25").Value = '  -11.79%  '
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("D27").Value = '3.24'
$ws.Range("E27").Value = '  -12.86%  '
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").Value = '8.06'
$ws.Range("E28").Value = '  -10.22%  '
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").Value = '29.38'
$ws.Range("E29").Value = '  -12.04%  '
$ws.Range("E30").Value = '  -16.20%  '
$ws.Range("E31").Value = '  -11.21%  '
$ws.Range("E32").Value = '  -11.05%  '
$ws.Range("D33").Value = '546.68'
$ws.Range("E33").Value = '  -11.84%  '
$ws.Range("D34").Value = '6.65'
$ws.Range("E34").Value = '  -17.60%  '
$ws.Range("D35").Value = '5.74'
$ws.Range("E35").Value = '  -15.29%  '
$ws.Range("E36").Value = '  +0.09%  '
$ws.Range("D37").Value = '0.0458'
$ws.Range("E37").Value = '  -6.59%  '
$ws.Range("D38").Value = '53.18'
$ws.Range("E38").Value = '  -6.10%  '
$ws.Range("D39").Value = '0.0861'
$ws.Range("E39").Value = '  -13.21%  '
$ws.Range("B40").Value = 'Cosmos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D40").Value = '9.11'
$ws.Range("E40").Value = '  -15.75%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").Value = '0.128'
$ws.Range("E41").Value = '  -9.52%  '
$ws.Range("D42").Value = '2.74'
$ws.Range("E42").Value = '  -19.48%  '
$ws.Range("D43").Value = '2.937.09'
$ws.Range("E43").Value = '  -11.69%  '
$ws.Range("B44").Value = 'PEPE'
$ws.Range("C44").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D44").Value = '0.0₃0591'
$ws.Range("E44").Value = '  -17.97%  '
$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D45").Value = '0.263'
$ws.Range("E45").Value = '  -14.91%  '
$ws.Range("D46").Value = '26.95'
$ws.Range("E46").Value = '  -14.97%  '
$ws.Range("D47").Value = '2.16'
$ws.Range("E47").Value = '  -13.37%  '
$ws.Range("B49").Value = 'ThetaToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D49").Value = '2.34'
$ws.Range("E49").Value = '  -19.97%  '
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").Value = '126.71'
$ws.Range("E50").Value = '  -5.37%  '
$ws.Range("E51").Value = '  -12.62%  '
